$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 16 de Septiembre de 2020 a las 17:09"

# Row 4
$ws.Cells.Item(4, 2).Value = 6791999
$ws.Cells.Item(4, 3).Value = 3852
$ws.Cells.Item(4, 4).Value = 4069609
$ws.Cells.Item(4, 5).Value = 2521925
$ws.Cells.Item(4, 7).Value = 268
$ws.Cells.Item(4, 8).Value = 200465

# Row 5
$ws.Cells.Item(5, 2).Value = 5060818
$ws.Cells.Item(5, 3).Value = 42784
$ws.Cells.Item(5, 4).Value = 3976413
$ws.Cells.Item(5, 5).Value = 1001901
$ws.Cells.Item(5, 7).Value = 413
$ws.Cells.Item(5, 8).Value = 82504

# Row 6
$ws.Cells.Item(6, 2).Value = 4384860
$ws.Cells.Item(6, 3).Value = 561
$ws.Cells.Item(6, 5).Value = 580515
$ws.Cells.Item(6, 7).Value = 10
$ws.Cells.Item(6, 8).Value = 133217

# Row 25
$ws.Cells.Item(25, 2).Value = 265964
$ws.Cells.Item(25, 3).Value = 1120
$ws.Cells.Item(25, 5).Value = 17417
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(25, 8).Value = 9447

# Row 29
$ws.Cells.Item(29, 2).Value = 139118
$ws.Cells.Item(29, 3).Value = 315
$ws.Cells.Item(29, 5).Value = 8088
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = 9190

# Row 31
$ws.Cells.Item(31, 2).Value = 122449
$ws.Cells.Item(31, 3).Value = 235
$ws.Cells.Item(31, 4).Value = 119400
$ws.Cells.Item(31, 5).Value = 2841

# Row 51
$ws.Cells.Item(51, 1).Value = "Portugal"
$ws.Cells.Item(51, 2).Value = 65626
$ws.Cells.Item(51, 3).Value = 605
$ws.Cells.Item(51, 4).Value = 44528
$ws.Cells.Item(51, 5).Value = 19220
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 1878

# Row 52
$ws.Cells.Item(52, 1).Value = "Etiopia"
$ws.Cells.Item(52, 2).Value = 65486
$ws.Cells.Item(52, 4).Value = 25988
$ws.Cells.Item(52, 5).Value = 38463
$ws.Cells.Item(52, 8).Value = 1035

# Row 61
$ws.Cells.Item(61, 5).Value = 6326
$ws.Cells.Item(61, 7).Value = 11
$ws.Cells.Item(61, 8).Value = 2039

# Row 65
$ws.Cells.Item(65, 2).Value = 44361
$ws.Cells.Item(65, 3).Value = 627
$ws.Cells.Item(65, 4).Value = 32732
$ws.Cells.Item(65, 5).Value = 10470
$ws.Cells.Item(65, 7).Value = 16
$ws.Cells.Item(65, 8).Value = 1159

# Row 75
$ws.Cells.Item(75, 2).Value = 27163
$ws.Cells.Item(75, 3).Value = 75
$ws.Cells.Item(75, 4).Value = 19960
$ws.Cells.Item(75, 5).Value = 6407

# Row 85
$ws.Cells.Item(85, 2).Value = 16088
$ws.Cells.Item(85, 3).Value = 163
$ws.Cells.Item(85, 4).Value = 13550
$ws.Cells.Item(85, 5).Value = 1870
$ws.Cells.Item(85, 7).Value = 7
$ws.Cells.Item(85, 8).Value = 668

# Row 93
$ws.Cells.Item(93, 2).Value = 12431
$ws.Cells.Item(93, 3).Value = 38
$ws.Cells.Item(93, 5).Value = 1795

# Row 94
$ws.Cells.Item(94, 2).Value = 11816
$ws.Cells.Item(94, 3).Value = 144
$ws.Cells.Item(94, 4).Value = 6733
$ws.Cells.Item(94, 5).Value = 4740
$ws.Cells.Item(94, 7).Value = 3
$ws.Cells.Item(94, 8).Value = 343

# Row 98
$ws.Cells.Item(98, 2).Value = 9964
$ws.Cells.Item(98, 3).Value = 63
$ws.Cells.Item(98, 4).Value = 7502
$ws.Cells.Item(98, 5).Value = 2354
$ws.Cells.Item(98, 7).Value = 2
$ws.Cells.Item(98, 8).Value = 108

# Row 141
$ws.Cells.Item(141, 1).Value = "Reunion"
$ws.Cells.Item(141, 2).Value = 3002
$ws.Cells.Item(141, 3).Value = 100
$ws.Cells.Item(141, 4).Value = 1313
$ws.Cells.Item(141, 5).Value = 1674
$ws.Cells.Item(141, 8).Value = 15

# Row 142
$ws.Cells.Item(142, 1).Value = "Mali"
$ws.Cells.Item(142, 2).Value = 2966
$ws.Cells.Item(142, 3).Value = 26
$ws.Cells.Item(142, 4).Value = 2311
$ws.Cells.Item(142, 5).Value = 527
$ws.Cells.Item(142, 8).Value = 128

# Row 168
$ws.Cells.Item(168, 4).Value = 936
$ws.Cells.Item(168, 5).Value = 92

# Row 204
$ws.Cells.Item(204, 1).Value = "Timor Oriental"

# Row 205
$ws.Cells.Item(205, 1).Value = "Santa Lucia"

# Row 214
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

# Row 215
$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1
